$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "Occurrence" column header to "Dates Used"
$ws.Cells.Item(1, 5).Value = "Dates Used"

# Replace the old occurrence/week-number codes in column E with the
# corresponding human-readable date ranges used for each row.
$ws.Cells.Item(2, 5).Value  = "Feb 3 - Feb 6, Feb 17 - Feb 19, March 31 - April 2"
$ws.Cells.Item(3, 5).Value  = "Feb 3 - Feb 6, Feb 17 - Feb 19, March 31 - April 2"
$ws.Cells.Item(4, 5).Value  = "Feb 3 - Feb 6, Feb 17 - Feb 19, March 31 - April 2"
$ws.Cells.Item(5, 5).Value  = "March 31 - April 2"
$ws.Cells.Item(6, 5).Value  = "Feb 17 - Feb 19, March 31 - April 2"
$ws.Cells.Item(7, 5).Value  = "Feb 17 - Feb 19, March 31 - April 2"
$ws.Cells.Item(8, 5).Value  = "Feb 3 - Feb 6, Feb 17 - Feb 19, March 31 - April 2"
$ws.Cells.Item(9, 5).Value  = "Feb 3 - Feb 6, Feb 17 - Feb 19, March 31 - April 2"
$ws.Cells.Item(10, 5).Value = "Feb 3 - Feb 6, Feb 17 - Feb 19, March 31 - April 2"
$ws.Cells.Item(11, 5).Value = "March 31 - April 2"
$ws.Cells.Item(12, 5).Value = "Feb 17 - Feb 19, March 31 - April 2"
$ws.Cells.Item(13, 5).Value = "Feb 17 - Feb 19, March 31 - April 2"
$ws.Cells.Item(14, 5).Value = "March 31 - April 2, May 12 - May 14"
$ws.Cells.Item(15, 5).Value = "March 31 - April 2"
$ws.Cells.Item(16, 5).Value = "March 31 - April 2, May 12 - May 14"
$ws.Cells.Item(17, 5).Value = "March 31 - April 2"
$ws.Cells.Item(18, 5).Value = "March 31 - April 2, May 12 - May 14"
$ws.Cells.Item(19, 5).Value = "March 31 - April 2, May 12 - May 14"
$ws.Cells.Item(20, 5).Value = "Feb 3 - Feb 6, Feb 17 - Feb 19, March 31 - April 2"
$ws.Cells.Item(21, 5).Value = "Feb 3 - Feb 6, Feb 17 - Feb 19, March 31 - April 2"
$ws.Cells.Item(22, 5).Value = "Feb 3 - Feb 6, Feb 17 - Feb 19, March 31 - April 2"
$ws.Cells.Item(23, 5).Value = "Feb 17 - Feb 19, March 31 - April 2"
$ws.Cells.Item(24, 5).Value = "Feb 17 - Feb 19, March 31 - April 2"
$ws.Cells.Item(25, 5).Value = "March 31 - April 2, May 12 - May 14"
$ws.Cells.Item(26, 5).Value = "March 31 - April 2"
$ws.Cells.Item(27, 5).Value = "March 31 - April 2, May 12 - May 14"
$ws.Cells.Item(28, 5).Value = "March 31 - April 2, May 12 - May 14"
$ws.Cells.Item(29, 5).Value = "March 31 - April 2, May 12 - May 14"

# Match the author's on-screen selection state after the edit.
$ws.Range("E1").Select()
